# Update the cryptocurrency price/volume table (rows 2-51) to the
# latest scraped values. Two coin pairs also swapped rank order
# (Maker/MXToken rows 37-38, BabyDogeCoin/Aptos/Algorand/EnergySwap
# rows 46-50), so Coin/Link/Price/Volume are rewritten wholesale per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the literal text into the cell (so strings that look like
    # numbers, e.g. "0.9990" or "1.0000", keep their exact formatting
    # instead of Excel silently parsing them into numeric values),
    # then drop the temporary Text number-format again so the cell
    # is left without any explicit style, matching the sheet original.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# IsNumericLooking marks the Price values that Excel would otherwise
# silently reinterpret as a plain number (dropping trailing zeroes,
# e.g. "0.9990" -> 0.999); the Volume(1h) column never needs this
# since its "  +x.xx%  " text never parses as a number.
$rows = @(
    @{ Row = 2; B = 'Bitcoin'; C = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D = '29.401.47'; DNumeric = $false; E = '  +0.11%  '; ENumeric = $false },
    @{ Row = 3; B = 'Ethereum'; C = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D = '1.846.43'; DNumeric = $false; E = '  +0.21%  '; ENumeric = $false },
    @{ Row = 4; B = 'TetherUSD'; C = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D = '0.9990'; DNumeric = $true; E = '  +0.00%  '; ENumeric = $false },
    @{ Row = 5; B = 'BNB'; C = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D = '240.50'; DNumeric = $true; E = '  +0.72%  '; ENumeric = $false },
    @{ Row = 6; B = 'XRP'; C = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D = '0.6265'; DNumeric = $true; E = '  -0.65%  '; ENumeric = $false },
    @{ Row = 7; B = 'USDC'; C = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D = '0.9998'; DNumeric = $true; E = '  -0.03%  '; ENumeric = $false },
    @{ Row = 8; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '0.07657'; DNumeric = $true; E = '  +1.73%  '; ENumeric = $false },
    @{ Row = 9; B = 'Cardano'; C = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D = '0.2915'; DNumeric = $true; E = '  -0.35%  '; ENumeric = $false },
    @{ Row = 10; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '24.81'; DNumeric = $true; E = '  +1.83%  '; ENumeric = $false },
    @{ Row = 11; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.07741'; DNumeric = $true; E = '  +0.43%  '; ENumeric = $false },
    @{ Row = 12; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.847.69'; DNumeric = $false; E = '  -0.50%  '; ENumeric = $false },
    @{ Row = 13; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '5.026'; DNumeric = $true; E = '  +0.71%  '; ENumeric = $false },
    @{ Row = 14; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '0.6799'; DNumeric = $true; E = '  +0.25%  '; ENumeric = $false },
    @{ Row = 15; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.00001073'; DNumeric = $true; E = '  +4.44%  '; ENumeric = $false },
    @{ Row = 16; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '83.41'; DNumeric = $true; E = '  +0.65%  '; ENumeric = $false },
    @{ Row = 17; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '6.167'; DNumeric = $true; E = '  +0.31%  '; ENumeric = $false },
    @{ Row = 18; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '29.425.91'; DNumeric = $false; E = '  +0.05%  '; ENumeric = $false },
    @{ Row = 19; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '228.02'; DNumeric = $true; E = '  +0.20%  '; ENumeric = $false },
    @{ Row = 20; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '12.38'; DNumeric = $true; E = '  -0.29%  '; ENumeric = $false },
    @{ Row = 21; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.0000'; DNumeric = $true; E = '  -0.04%  '; ENumeric = $false },
    @{ Row = 22; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '7.408'; DNumeric = $true; E = '  -0.43%  '; ENumeric = $false },
    @{ Row = 23; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.000'; DNumeric = $true; E = '  -0.06%  '; ENumeric = $false },
    @{ Row = 24; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '157.87'; DNumeric = $true; E = '  +0.60%  '; ENumeric = $false },
    @{ Row = 25; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.1371'; DNumeric = $true; E = '  -1.31%  '; ENumeric = $false },
    @{ Row = 26; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '8.395'; DNumeric = $true; E = '  +0.49%  '; ENumeric = $false },
    @{ Row = 27; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '17.68'; DNumeric = $true; E = '  +0.48%  '; ENumeric = $false },
    @{ Row = 28; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '1.354'; DNumeric = $true; E = '  +6.08%  '; ENumeric = $false },
    @{ Row = 29; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '1.461'; DNumeric = $true; E = '  +0.39%  '; ENumeric = $false },
    @{ Row = 30; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.05656'; DNumeric = $true; E = '  +0.62%  '; ENumeric = $false },
    @{ Row = 31; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '4.115'; DNumeric = $true; E = '  +0.28%  '; ENumeric = $false },
    @{ Row = 32; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '4.025'; DNumeric = $true; E = '  +0.17%  '; ENumeric = $false },
    @{ Row = 33; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '1.840'; DNumeric = $true; E = '  +0.58%  '; ENumeric = $false },
    @{ Row = 34; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.160'; DNumeric = $true; E = '  +0.44%  '; ENumeric = $false },
    @{ Row = 35; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '0.7007'; DNumeric = $true; E = '  -0.77%  '; ENumeric = $false },
    @{ Row = 36; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '2.587'; DNumeric = $true; E = '  -0.06%  '; ENumeric = $false },
    @{ Row = 37; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '2.761'; DNumeric = $true; E = '  -0.11%  '; ENumeric = $false },
    @{ Row = 38; B = 'Maker'; C = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D = '1.222.51'; DNumeric = $false; E = '  -1.48%  '; ENumeric = $false },
    @{ Row = 39; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.01787'; DNumeric = $true; E = '  -0.92%  '; ENumeric = $false },
    @{ Row = 40; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '6.540'; DNumeric = $true; E = '  +3.60%  '; ENumeric = $false },
    @{ Row = 41; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '0.9047'; DNumeric = $true; E = '  +0.54%  '; ENumeric = $false },
    @{ Row = 42; B = 'PaxDollar'; C = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D = '1.000'; DNumeric = $true; E = '  +0.07%  '; ENumeric = $false },
    @{ Row = 43; B = 'RocketPoolETH'; C = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'; D = '2.006.66'; DNumeric = $false; E = '  -1.07%  '; ENumeric = $false },
    @{ Row = 44; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '101.71'; DNumeric = $true; E = '  -0.19%  '; ENumeric = $false },
    @{ Row = 45; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '65.91'; DNumeric = $true; E = '  +0.50%  '; ENumeric = $false },
    @{ Row = 46; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '7.147'; DNumeric = $true; E = '  +1.25%  '; ENumeric = $false },
    @{ Row = 47; B = 'BabyDogeCoin'; C = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D = '0.00000000119'; DNumeric = $true; E = '  -0.47%  '; ENumeric = $false },
    @{ Row = 48; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.4011'; DNumeric = $true; E = '  +0.31%  '; ENumeric = $false },
    @{ Row = 49; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '9.022'; DNumeric = $true; E = '  +1.76%  '; ENumeric = $false },
    @{ Row = 50; B = 'Algorand'; C = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D = '0.1153'; DNumeric = $true; E = '  +3.11%  '; ENumeric = $false },
    @{ Row = 51; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '1.670'; DNumeric = $true; E = '  +0.31%  '; ENumeric = $false }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C

    $dCell = $ws.Cells.Item($r.Row, 4)
    if ($r.DNumeric) {
        Set-TextValue $dCell $r.D
    } else {
        $dCell.Value = $r.D
    }

    $eCell = $ws.Cells.Item($r.Row, 5)
    if ($r.ENumeric) {
        Set-TextValue $eCell $r.E
    } else {
        $eCell.Value = $r.E
    }
}
